$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 33672, "Maria Luísa Azevedo", "Juridico", "Doenca", 1, 45096, 5096.6),
    @(3, 39724, "Dra. Maria Helena Fogaça", "Marketing", "Viagem de negocios", 6, 45101, 6967.85),
    @(4, 41566, "Daniel Moraes", "Engenharia", "Viagem de negocios", 8, 45104, 5924.76),
    @(5, 46584, "Dom Leão", "Atendimento ao Cliente", "Doenca", 4, 45082, 9149.139999999999),
    @(6, 3437, "Lavínia Ramos", "Recursos Humanos", "Problemas pessoais", 4, 45088, 9181.6),
    @(7, 44811, "Maria Flor Lima", "Atendimento ao Cliente", "Viagem de negocios", 7, 45096, 7749.52),
    @(8, 10653, "Srta. Lorena Mendes", "Atendimento ao Cliente", "Viagem de negocios", 2, 45087, 8928.18),
    @(9, 62258, "Carlos Eduardo Campos", "Vendas", "Outros", 4, 45105, 8445.75),
    @(10, 6761, "Dr. Cauê Oliveira", "P&D", "Problemas pessoais", 5, 45099, 5006.39),
    @(11, 393, "Ana Lívia da Mata", "Juridico", "Viagem de negocios", 1, 45089, 8990.379999999999)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
